# docs/diagrams/StorageComponentClassDiagram.pptx
# - Update interface AddressBookStorage to SchedulerStorage
# - Update class XmlAddressBookStorage to XmlSchedulerStorage
# - Update class XmlSerializableAddressBook to XmlSerializableScheduler
# - Update class XmlAdaptedPerson to XmlAdaptedEvent

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# Shape id=2 "Rectangle 8" (interface box):
#   "<<interface>>" + (line break) + "AddressBookStorage"
#   -> "<<interface>>" + (line break) + "SchedulerStorage"
$shp = Get-ShapeById $s 2
$tr = $shp.TextFrame.TextRange
$tr.Characters(15, 18).Text = "SchedulerStorage"

# Shape id=50 "Rectangle 8" (class box):
#   "XmlAddressBook" + (line break) + "Storage"
#   -> "XmlScheduler" + (line break) + "Storage"
$shp = Get-ShapeById $s 50
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 14).Text = "XmlScheduler"

# Shape id=66 "Rectangle 8" (class box):
#   "XmlSerializable" + (line break) + "AddressBook"
#   -> "XmlSerializable" + (line break) + "Scheduler"
$shp = Get-ShapeById $s 66
$tr = $shp.TextFrame.TextRange
$tr.Characters(17, 11).Text = "Scheduler"

# Shape id=74 "Rectangle 8" (class box):
#   "XmlAdaptedPerson" -> "XmlAdaptedEvent"
$shp = Get-ShapeById $s 74
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 16).Text = "XmlAdaptedEvent"
